# Apply "add colorful rtx3070 ad" changes to GraphicsCardSpecificationSheet.xlsx
#
# Summary of changes:
#  1. Sheet "RTX3080 10G", cell A7: fix a copy/paste typo
#       "Colorful iGame GeForce RTX3090 Advanced OC"
#       -> "Colorful iGame GeForce RTX3080 Advanced OC"
#  2. Sheet "RTX3070 8G": clarify the MOS text in rows 2-4 (F column) to
#     distinguish GPU vs Memory power stages.
#  3. Sheet "RTX3070 8G": populate the previously-empty row 6 with the
#     new "Colorful iGame GeForce RTX3070 Advanced OC" card data.
#  4. Row-height / wrap adjustments that follow from the above content
#     changes, and refreshed cell selections.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) RTX3080 10G : fix A7 typo
# ---------------------------------------------------------------------
$ws80 = $wb.Worksheets.Item("RTX3080 10G")
$ws80.Cells.Item(7, 1).Value = "Colorful iGame GeForce RTX3080 Advanced OC"

# ---------------------------------------------------------------------
# 2) + 3) RTX3070 8G updates
# ---------------------------------------------------------------------
$ws70 = $wb.Worksheets.Item("RTX3070 8G")

# -- row 2 : MOS (F2) now explicitly calls out GPU vs Memory stages,
#    and PWR Connectors (G2) switches from a generic "10-phase" mix
#    note to the actual connector description.
$ws70.Cells.Item(2, 6).Value = "AOS  `nAOZ5311`n(50A DrMOS,GPU)`nSinopower SM7342EKKP`n(Memory)"
$ws70.Cells.Item(2, 7).Value = "1*12PIN`n(only 6 pins are valid)`n"

# -- row 3 : MOS (F3)
$ws70.Cells.Item(3, 6).Value = "AOS  `nAOZ5311`n(50A DrMOS,GPU)`nOnSemi `nON3102`nON3106`n(Memory)"

# -- row 4 : MOS (F4)
$ws70.Cells.Item(4, 6).Value = "TI NexFET `nCSD95481RWJ`n(60A DrMOS,GPU)`nOnsemi NCP303151`n(Memory)"

# -- row 6 : brand new card entry - Colorful iGame GeForce RTX3070 Advanced OC
$ws70.Cells.Item(6, 1).Value = "Colorful iGame GeForce RTX3070 Advanced OC"
$ws70.Cells.Item(6, 2).Value = "10-phase"
$ws70.Cells.Item(6, 3).Value = "3-phase"
$ws70.Cells.Item(6, 4).Value = "UPI uP9512R"
$ws70.Cells.Item(6, 5).Value = "?"
$ws70.Cells.Item(6, 6).Value = "AOS  `nAOZ5311`n(50A DrMOS,GPU)`nSinopower`nSM4364A`nSM4373`n(Memory)`n"
$ws70.Cells.Item(6, 7).Value = "2*8PIN"
$ws70.Cells.Item(6, 8).Value = "5*8mm heat pipe with vapor chamber`nmetal backplate with thermal pads"
$ws70.Cells.Item(6, 9).Value = "2*9cm+1*8cm`n3000rpm"
$ws70.Cells.Item(6, 10).Value = "1500/`n1815/`n1995?MHZ"
$ws70.Cells.Item(6, 11).Value = "1750MHZ"
$ws70.Cells.Item(6, 12).Value = "270/290W"
$ws70.Cells.Item(6, 13).Value = "66°C/2258rpm`n(expreview,25°C)`n70°C/2280rpm`n(chiphell,20°C)"
$ws70.Cells.Item(6, 14).Value = "316*131*53`nmm"
$ws70.Cells.Item(6, 15).Value = "1*HDMI`n3*DP"
$ws70.Cells.Item(6, 16).Value = "www.chiphell.com/portal.php?mod=view&aid=24447&page=5`nwww.expreview.com/76459.html"

# Column G lost its wrap-text styling now that it holds a short value
# (matches the style used for the same column in the other data rows).
$ws70.Cells.Item(6, 7).WrapText = $false

# Row heights follow the new (taller) wrapped content.
$ws70.Rows.Item(2).RowHeight = 99.75
$ws70.Rows.Item(3).RowHeight = 114
$ws70.Rows.Item(6).RowHeight = 128.25

# ---------------------------------------------------------------------
# Refresh the active cell / selection on each sheet, as recorded by Excel
# ---------------------------------------------------------------------
$ws3090 = $wb.Worksheets.Item("RTX3090 24G")
[void]$ws3090.Activate()
[void]$ws3090.Range("A7").Select()

[void]$ws80.Activate()
[void]$ws80.Range("P7").Select()

[void]$ws70.Activate()
[void]$ws70.Range("P7").Select()
